$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder player blocks (rows 11-28, 32-43) alphabetically by player name,
# updating Player (A), PD (C), Comb (D), Solo (E), Ast (F) - B (Season Group) unchanged.

$ws.Range("A11").Value = "Harrison Smith"
$ws.Range("C11").Value = 9.333333333333334
$ws.Range("D11").Value = 96
$ws.Range("E11").Value = 67.33333333333333
$ws.Range("F11").Value = 28.66666666666667

$ws.Range("A12").Value = "Harrison Smith"
$ws.Range("C12").Value = 7.666666666666667
$ws.Range("D12").Value = 88.33333333333333
$ws.Range("E12").Value = 60.66666666666666
$ws.Range("F12").Value = 27.66666666666667

$ws.Range("A13").Value = "Harrison Smith"
$ws.Range("C13").Value = -1.666666666666667
$ws.Range("D13").Value = -7.666666666666671
$ws.Range("E13").Value = -6.666666666666664
$ws.Range("F13").Value = -1

$ws.Range("A14").Value = "Jalen Mills"
$ws.Range("C14").Value = 5.666666666666667
$ws.Range("D14").Value = 54
$ws.Range("E14").Value = 41
$ws.Range("F14").Value = 13

$ws.Range("A15").Value = "Jalen Mills"
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 40
$ws.Range("E15").Value = 23
$ws.Range("F15").Value = 17

$ws.Range("A16").Value = "Jalen Mills"
$ws.Range("C16").Value = -1.666666666666667
$ws.Range("D16").Value = -14
$ws.Range("E16").Value = -18
$ws.Range("F16").Value = 4

$ws.Range("A17").Value = "Jimmie Ward"
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 71.66666666666667
$ws.Range("E17").Value = 49
$ws.Range("F17").Value = 22.66666666666667

$ws.Range("A18").Value = "Jimmie Ward"
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 49.33333333333334
$ws.Range("E18").Value = 34.33333333333334
$ws.Range("F18").Value = 15

$ws.Range("A19").Value = "Jimmie Ward"
$ws.Range("C19").Value = -2
$ws.Range("D19").Value = -22.33333333333334
$ws.Range("E19").Value = -14.66666666666666
$ws.Range("F19").Value = -7.666666666666668

$ws.Range("A20").Value = "Julian Love"
$ws.Range("C20").Value = 4.333333333333333
$ws.Range("D20").Value = 55.66666666666666
$ws.Range("E20").Value = 37.66666666666666
$ws.Range("F20").Value = 18

$ws.Range("A21").Value = "Julian Love"
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 118.6666666666667
$ws.Range("E21").Value = 81
$ws.Range("F21").Value = 37.66666666666666

$ws.Range("A22").Value = "Julian Love"
$ws.Range("C22").Value = 4.666666666666667
$ws.Range("D22").Value = 63.00000000000001
$ws.Range("E22").Value = 43.33333333333334
$ws.Range("F22").Value = 19.66666666666666

$ws.Range("A23").Value = "Justin Reid"
$ws.Range("C23").Value = 4.333333333333333
$ws.Range("D23").Value = 75.66666666666667
$ws.Range("E23").Value = 53.33333333333334
$ws.Range("F23").Value = 22.33333333333333

$ws.Range("A24").Value = "Justin Reid"
$ws.Range("C24").Value = 7.666666666666667
$ws.Range("D24").Value = 88.33333333333333
$ws.Range("E24").Value = 64.66666666666667
$ws.Range("F24").Value = 23.66666666666667

$ws.Range("A25").Value = "Justin Reid"
$ws.Range("C25").Value = 3.333333333333334
$ws.Range("D25").Value = 12.66666666666666
$ws.Range("E25").Value = 11.33333333333334
$ws.Range("F25").Value = 1.333333333333336

$ws.Range("A26").Value = "Lonnie Johnson"
$ws.Range("C26").Value = 4.333333333333333
$ws.Range("D26").Value = 57.33333333333334
$ws.Range("E26").Value = 42.66666666666666
$ws.Range("F26").Value = 14.66666666666667

$ws.Range("A27").Value = "Lonnie Johnson"
$ws.Range("C27").Value = 0.6666666666666666
$ws.Range("D27").Value = 9
$ws.Range("E27").Value = 7.333333333333333
$ws.Range("F27").Value = 1.666666666666667

$ws.Range("A28").Value = "Lonnie Johnson"
$ws.Range("C28").Value = -3.666666666666667
$ws.Range("D28").Value = -48.33333333333334
$ws.Range("E28").Value = -35.33333333333333
$ws.Range("F28").Value = -13

$ws.Range("A32").Value = "Minkah Fitzpatrick"
$ws.Range("C32").Value = 8
$ws.Range("D32").Value = 83
$ws.Range("E32").Value = 57.77777777777777
$ws.Range("F32").Value = 25.22222222222222

$ws.Range("A33").Value = "Minkah Fitzpatrick"
$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 85.33333333333333
$ws.Range("E33").Value = 53.66666666666666
$ws.Range("F33").Value = 31.66666666666667

$ws.Range("A34").Value = "Minkah Fitzpatrick"
$ws.Range("C34").Value = -2
$ws.Range("D34").Value = 2.333333333333329
$ws.Range("E34").Value = -4.111111111111107
$ws.Range("F34").Value = 6.444444444444443

$ws.Range("A35").Value = "Rodney McLeod"
$ws.Range("C35").Value = 5.666666666666667
$ws.Range("D35").Value = 66.66666666666667
$ws.Range("E35").Value = 40.33333333333334
$ws.Range("F35").Value = 26.33333333333333

$ws.Range("A36").Value = "Rodney McLeod"
$ws.Range("C36").Value = 4.666666666666667
$ws.Range("D36").Value = 54.66666666666666
$ws.Range("E36").Value = 34.33333333333334
$ws.Range("F36").Value = 20.33333333333333

$ws.Range("A37").Value = "Rodney McLeod"
$ws.Range("C37").Value = -1
$ws.Range("D37").Value = -12.00000000000001
$ws.Range("E37").Value = -6
$ws.Range("F37").Value = -6

$ws.Range("A38").Value = "Ronnie Harrison"
$ws.Range("C38").Value = 6.333333333333333
$ws.Range("D38").Value = 55.66666666666666
$ws.Range("E38").Value = 36
$ws.Range("F38").Value = 19.66666666666667

$ws.Range("A39").Value = "Ronnie Harrison"
$ws.Range("C39").Value = 1
$ws.Range("D39").Value = 15.33333333333333
$ws.Range("E39").Value = 9.666666666666666
$ws.Range("F39").Value = 5.666666666666667

$ws.Range("A40").Value = "Ronnie Harrison"
$ws.Range("C40").Value = -5.333333333333333
$ws.Range("D40").Value = -40.33333333333333
$ws.Range("E40").Value = -26.33333333333334
$ws.Range("F40").Value = -14

$ws.Range("A41").Value = "Tyrann Mathieu"
$ws.Range("C41").Value = 9
$ws.Range("D41").Value = 71
$ws.Range("E41").Value = 57
$ws.Range("F41").Value = 14

$ws.Range("A42").Value = "Tyrann Mathieu"
$ws.Range("C42").Value = 8
$ws.Range("D42").Value = 76
$ws.Range("E42").Value = 53.66666666666666
$ws.Range("F42").Value = 22.33333333333333

$ws.Range("A43").Value = "Tyrann Mathieu"
$ws.Range("C43").Value = -1
$ws.Range("D43").Value = 5
$ws.Range("E43").Value = -3.333333333333336
$ws.Range("F43").Value = 8.333333333333332
